$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.934.17"
$ws.Range("E2").Value = "  -0.14%  "

$ws.Range("D3").Value = "2.116.89"
$ws.Range("E3").Value = "  +0.88%  "

$ws.Range("D4").Value = "'1.009"
$ws.Range("E4").Value = "  +0.55%  "

$ws.Range("D5").Value = "'348.08"
$ws.Range("E5").Value = "  +0.36%  "

$ws.Range("D6").Value = "'1.007"
$ws.Range("E6").Value = "  +0.48%  "

$ws.Range("D7").Value = "'0.5211"
$ws.Range("E7").Value = "  +0.97%  "

$ws.Range("D8").Value = "'0.4452"
$ws.Range("E8").Value = "  +0.66%  "

$ws.Range("D9").Value = "'54.60"
$ws.Range("E9").Value = "  +4.50%  "

$ws.Range("D10").Value = "'0.09359"
$ws.Range("E10").Value = "  -0.92%  "

$ws.Range("E11").Value = "  +0.91%  "

$ws.Range("D12").Value = "'25.18"
$ws.Range("E12").Value = "  -0.26%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "2.111.58"
$ws.Range("E13").Value = "  +0.44%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'8.427"
$ws.Range("E14").Value = "  +3.16%  "

$ws.Range("D15").Value = "'6.868"
$ws.Range("E15").Value = "  +1.77%  "

$ws.Range("D16").Value = "'102.23"
$ws.Range("E16").Value = "  +3.10%  "

$ws.Range("D17").Value = "'0.00001164"
$ws.Range("E17").Value = "  -0.47%  "

$ws.Range("D18").Value = "'1.008"
$ws.Range("E18").Value = "  +0.41%  "

$ws.Range("D19").Value = "'21.52"
$ws.Range("E19").Value = "  +4.52%  "

$ws.Range("D20").Value = "'0.06686"
$ws.Range("E20").Value = "  +0.18%  "

$ws.Range("E21").Value = "  +1.28%  "

$ws.Range("E22").Value = "  +0.48%  "

$ws.Range("D23").Value = "29.958.06"
$ws.Range("E23").Value = "  -0.41%  "

$ws.Range("D24").Value = "'12.75"
$ws.Range("E24").Value = "  +1.11%  "

$ws.Range("D25").Value = "'2.326"
$ws.Range("E25").Value = "  -0.22%  "

$ws.Range("D26").Value = "2.348.36"
$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("E27").Value = "  +0.69%  "

$ws.Range("D28").Value = "'2.552"
$ws.Range("E28").Value = "  -0.15%  "

$ws.Range("D29").Value = "'162.37"
$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("D30").Value = "'133.95"
$ws.Range("E30").Value = "  +0.42%  "

$ws.Range("D31").Value = "'1.157"
$ws.Range("E31").Value = "  -1.10%  "

$ws.Range("D32").Value = "'1.791"
$ws.Range("E32").Value = "  +9.24%  "

$ws.Range("E33").Value = "  -0.43%  "

$ws.Range("D34").Value = "'6.242"
$ws.Range("E34").Value = "  +0.32%  "

$ws.Range("D35").Value = "'3.968"

$ws.Range("D36").Value = "'6.507"
$ws.Range("E36").Value = "  +4.87%  "

$ws.Range("D37").Value = "'10.82"
$ws.Range("E37").Value = "  +7.51%  "

$ws.Range("E38").Value = "  +1.47%  "

$ws.Range("D39").Value = "'0.06859"
$ws.Range("E39").Value = "  +1.15%  "

$ws.Range("D40").Value = "'0.7020"
$ws.Range("E40").Value = "  +1.01%  "

$ws.Range("E41").Value = "  +1.45%  "

$ws.Range("D42").Value = "'0.2247"
$ws.Range("E42").Value = "  -1.21%  "

$ws.Range("D43").Value = "'1.331"

$ws.Range("D44").Value = "'0.6827"
$ws.Range("E44").Value = "  +3.04%  "

$ws.Range("D45").Value = "'14.56"
$ws.Range("E45").Value = "  +2.77%  "

$ws.Range("D46").Value = "'2.353"
$ws.Range("E46").Value = "  +3.24%  "

$ws.Range("E47").Value = "  +0.47%  "

$ws.Range("D48").Value = "'3.637"
$ws.Range("E48").Value = "  +0.16%  "

$ws.Range("D49").Value = "'0.00000000355"
$ws.Range("E49").Value = "  +0.96%  "

$ws.Range("D50").Value = "'1.244"
$ws.Range("E50").Value = "  +6.91%  "

$ws.Range("D51").Value = "'1.223"
$ws.Range("E51").Value = "  +0.16%  "
